$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$genderQuery = "MATCH (s:study)<--(p:participant)
WHERE s.study_name in [""Detection of Colorectal Cancer Susceptibility Loci Using Genome-Wide Sequencing""]
OPTIONAL MATCH (p)<--(samp:sample)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN   
 coalesce(p.participant_id,'') as ``Participant ID``,
 coalesce(s.study_name, '') as ``Study Name``,
 coalesce(s.phs_accession,'') as ``Accession``,
 coalesce(p.gender,'') as ``Gender``,
 coalesce(apoc.text.join(samp, ','), '') as ``Samples``
 ORDER By p.participant_id LIMIT 100"

$tumorQuery = "MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in [""Detection of Colorectal Cancer Susceptibility Loci Using Genome-Wide Sequencing""]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as ``Sample ID``,
 coalesce(p.participant_id,'') as ``Participant ID``,
 coalesce(s.study_name, '') as ``Study Name``,
 coalesce(s.phs_accession,'') as ``Accession``,
coalesce(samp.sample_tumor_status,'') as ``Tumor``,
coalesce(samp.sample_type,'') as ``Analyte Type``
ORDER By samp.sample_id LIMIT 100"

# Row 2 = CasesTab -> participant/Gender query
$ws.Cells.Item(2, 2).Value = $genderQuery

# Row 3 = SamplesTab -> sample/Tumor query, fixed to use samp.sample_tumor_status
$ws.Cells.Item(3, 2).Value = $tumorQuery

# Restore the cursor/selection position recorded in the saved view state
$ws.Range("C11").Select()
